$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new row 6 cells (order matters for shared-string index allocation)
$ws.Range("B6").Value = "https://api.nasa.gov/mars-photos/api/v1/rovers/curiosity/photos?sol=1000&page=2&api_key=DEMO_KEY"
$ws.Range("C6").Value = "photos[1].rover.name"
$ws.Range("D6").Value = "Curiosity"
$ws.Range("A6").Value = "TestCase_005"

# Update existing JsonValueExpected cell for TestCase_004 (video -> image)
$ws.Range("D5").Value = "image"

# Copy E5's cell (keeps the quote-prefixed "200" text style) down into E6
$ws.Range("E5").Copy($ws.Range("E6"))

# Update the selected cell shown when the workbook is reopened
$ws.Range("D5").Select()
